$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Data" sheet: append 5 new daily observations after row 443 (A1:B443 ->
#    A1:B448), continuing the existing date/value series.
# ---------------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("Data")

$newRows = @(
    @{ row = 444; date = 45120; val = 1767.432 },
    @{ row = 445; date = 45121; val = 1740.777 },
    @{ row = 446; date = 45124; val = 1728.322 },
    @{ row = 447; date = 45125; val = 1716.862 },
    @{ row = 448; date = 45126; val = 1732.804 }
)

foreach ($r in $newRows) {
    # Copy the format (style) of the last existing row down into the new row
    # so the date column keeps its date-formatted style.
    $src = $wsData.Range("A443:B443")
    $dst = $wsData.Range("A" + $r.row + ":B" + $r.row)
    $src.Copy($dst)

    $wsData.Cells.Item($r.row, 1).Value = $r.date
    $wsData.Cells.Item($r.row, 2).Value = $r.val
}

# ---------------------------------------------------------------------------
# 2. "SeriesInfo" sheet: refresh the metadata fields that changed.
# ---------------------------------------------------------------------------
$wsInfo = $wb.Worksheets.Item("SeriesInfo")

function Set-PlainTextValue($cell, $text) {
    # Assigning an ISO-date-looking string straight to .Value causes Excel's
    # smart-parser to convert it into a date serial number. Going through a
    # text formula and then converting that formula to a static value keeps
    # the cell a genuine text value (matching the source file) without
    # picking up stray number formats/styles.
    $cell.Formula = '="' + $text + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163) # xlPasteValues
}

Set-PlainTextValue $wsInfo.Cells.Item(3, 2) "2023-07-20"
Set-PlainTextValue $wsInfo.Cells.Item(4, 2) "2023-07-20"
Set-PlainTextValue $wsInfo.Cells.Item(7, 2) "2023-07-19"
$wsInfo.Cells.Item(14, 2).Value = "2023-07-19 13:01:03-05"
